$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.240.01"
$ws.Range("E2").Value = "  -0.20%  "
Set-TextValue $ws.Range("D3") "3.404.99"
$ws.Range("E3").Value = "  +0.86%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "573.28"
$ws.Range("E5").Value = "  +0.00%  "
Set-TextValue $ws.Range("D6") "139.34"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue $ws.Range("D8") "3.403.48"
$ws.Range("E8").Value = "  +0.78%  "
Set-TextValue $ws.Range("D9") "0.471"
$ws.Range("E9").Value = "  -0.34%  "
Set-TextValue $ws.Range("D10") "7.71"
$ws.Range("E10").Value = "  +3.19%  "
Set-TextValue $ws.Range("D11") "0.123"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  -1.76%  "
Set-TextValue $ws.Range("D13") "3.980.14"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  -1.07%  "
Set-TextValue $ws.Range("D15") "26.67"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  -1.50%  "
Set-TextValue $ws.Range("D17") "3.396.47"
$ws.Range("E17").Value = "  +0.56%  "
Set-TextValue $ws.Range("D18") "61.264.45"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +1.35%  "
Set-TextValue $ws.Range("D20") "13.90"
$ws.Range("E20").Value = "  -0.77%  "
Set-TextValue $ws.Range("D21") "9.32"
$ws.Range("E21").Value = "  -0.31%  "
Set-TextValue $ws.Range("D22") "378.00"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D23") "0.554"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D24") "3.532.46"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +0.08%  "
Set-TextValue $ws.Range("D26") "71.35"
$ws.Range("E26").Value = "  -0.43%  "
Set-TextValue $ws.Range("D27") "0.0000123"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("E28").Value = "  +9.72%  "
$ws.Range("E29").Value = "  -7.38%  "
$ws.Range("E30").Value = "  -0.13%  "
Set-TextValue $ws.Range("D31") "7.44"
$ws.Range("E31").Value = "  -1.31%  "
Set-TextValue $ws.Range("D32") "8.18"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -0.01%  "
Set-TextValue $ws.Range("D35") "23.48"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "5.13"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D38") "6.90"
$ws.Range("E38").Value = "  +1.36%  "
Set-TextValue $ws.Range("D39") "166.34"
$ws.Range("E39").Value = "  +0.53%  "
Set-TextValue $ws.Range("D40") "0.0775"
$ws.Range("E40").Value = "  +0.38%  "
Set-TextValue $ws.Range("D41") "26.15"
$ws.Range("E41").Value = "  +7.27%  "
Set-TextValue $ws.Range("D42") "1.76"
$ws.Range("E42").Value = "  +1.52%  "
Set-TextValue $ws.Range("D43") "0.999"
$ws.Range("E43").Value = "  -0.12%  "
Set-TextValue $ws.Range("D44") "0.778"
$ws.Range("E44").Value = "  +0.68%  "
Set-TextValue $ws.Range("D45") "41.99"
$ws.Range("E45").Value = "  +0.88%  "
Set-TextValue $ws.Range("D46") "4.41"
$ws.Range("E46").Value = "  +0.14%  "
Set-TextValue $ws.Range("D47") "1.18"
$ws.Range("E47").Value = "  -2.46%  "
Set-TextValue $ws.Range("D48") "2.532.55"
$ws.Range("E48").Value = "  +7.10%  "
Set-TextValue $ws.Range("D49") "23.88"
$ws.Range("E49").Value = "  +5.64%  "
Set-TextValue $ws.Range("D50") "6.79"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  -0.05%  "
